$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.0843278244137764
$ws.Cells.Item(2, 2).Value = 0.979752242565155
$ws.Cells.Item(2, 3).Value = 0.02071941457688808
$ws.Cells.Item(2, 4).Value = 0.9967319965362549
$ws.Cells.Item(3, 1).Value = 0.01335260644555092
$ws.Cells.Item(3, 2).Value = 0.997761070728302
$ws.Cells.Item(3, 3).Value = 0.0151249598711729
$ws.Cells.Item(3, 4).Value = 0.9967319965362549
$ws.Cells.Item(4, 1).Value = 0.006928476504981518
$ws.Cells.Item(4, 2).Value = 0.9984667897224426
$ws.Cells.Item(4, 3).Value = 0.008190981112420559
$ws.Cells.Item(4, 4).Value = 0.9970291256904602
$ws.Cells.Item(5, 1).Value = 0.004671342670917511
$ws.Cells.Item(5, 2).Value = 0.9987831711769104
$ws.Cells.Item(5, 3).Value = 0.002630555303767323
$ws.Cells.Item(5, 4).Value = 0.998811662197113
$ws.Cells.Item(6, 1).Value = 0.002070237649604678
$ws.Cells.Item(6, 2).Value = 0.9994646310806274
$ws.Cells.Item(6, 3).Value = 0.001330101396888494
$ws.Cells.Item(6, 4).Value = 0.9994058012962341
$ws.Cells.Item(7, 1).Value = 0.002344866283237934
$ws.Cells.Item(7, 2).Value = 0.9994402527809143
$ws.Cells.Item(7, 3).Value = 0.003938700072467327
$ws.Cells.Item(7, 4).Value = 0.9994058012962341
$ws.Cells.Item(8, 1).Value = 0.001364344730973244
$ws.Cells.Item(8, 2).Value = 0.9996592998504639
$ws.Cells.Item(8, 3).Value = 0.001193636679090559
$ws.Cells.Item(8, 4).Value = 0.9991087317466736
$ws.Cells.Item(9, 1).Value = 0.001214924734085798
$ws.Cells.Item(9, 2).Value = 0.9996836185455322
$ws.Cells.Item(9, 3).Value = 0.001525252358987927
$ws.Cells.Item(9, 4).Value = 0.9991087317466736
$ws.Cells.Item(10, 1).Value = 0.001690750010311604
$ws.Cells.Item(10, 2).Value = 0.9997079372406006
$ws.Cells.Item(10, 3).Value = 0.0002496922970749438
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(11, 1).Value = 0.001630857354030013
$ws.Cells.Item(11, 2).Value = 0.9996349811553955
$ws.Cells.Item(11, 3).Value = 0.0003392604412510991
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(12, 1).Value = 0.0005658797454088926
$ws.Cells.Item(12, 2).Value = 0.9998783469200134
$ws.Cells.Item(12, 3).Value = 0.0003017795679625124
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(13, 1).Value = 0.000611230090726167
$ws.Cells.Item(13, 2).Value = 0.9997809529304504
$ws.Cells.Item(13, 3).Value = 0.00008377544145332649
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(14, 1).Value = 0.0006752714980393648
$ws.Cells.Item(14, 2).Value = 0.9998053312301636
$ws.Cells.Item(14, 3).Value = 0.0001135682323365472
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(15, 1).Value = 0.000718479510396719
$ws.Cells.Item(15, 2).Value = 0.9999026656150818
$ws.Cells.Item(15, 3).Value = 0.00001963444628927391
$ws.Cells.Item(15, 4).Value = 1
$ws.Cells.Item(16, 1).Value = 0.0002554840466473252
$ws.Cells.Item(16, 2).Value = 0.9999269843101501
$ws.Cells.Item(16, 3).Value = 0.00001591157342772931
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(17, 1).Value = 0.0007509626448154449
$ws.Cells.Item(17, 2).Value = 0.9998053312301636
$ws.Cells.Item(17, 3).Value = 0.0003228774294257164
$ws.Cells.Item(17, 4).Value = 0.9997029304504395
$ws.Cells.Item(18, 1).Value = 0.001732100965455174
$ws.Cells.Item(18, 2).Value = 0.9996836185455322
$ws.Cells.Item(18, 3).Value = 0.00003913978434866294
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(19, 1).Value = 0.0006746989674866199
$ws.Cells.Item(19, 2).Value = 0.9999026656150818
$ws.Cells.Item(19, 3).Value = 0.00007824574277037755
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(20, 1).Value = 0.0002382063103141263
$ws.Cells.Item(20, 2).Value = 0.9999513030052185
$ws.Cells.Item(20, 3).Value = 0.0005907302838750184
$ws.Cells.Item(20, 4).Value = 0.9997029304504395
$ws.Cells.Item(21, 1).Value = 0.00009719676017994061
$ws.Cells.Item(21, 2).Value = 0.9999756813049316
$ws.Cells.Item(21, 3).Value = 0.0004244929878041148
$ws.Cells.Item(21, 4).Value = 0.9997029304504395
$ws.Cells.Item(22, 1).Value = 0.0001249003398697823
$ws.Cells.Item(22, 2).Value = 0.9999756813049316
$ws.Cells.Item(22, 3).Value = 0.000006544911684613908
$ws.Cells.Item(22, 4).Value = 1
$ws.Cells.Item(23, 1).Value = 0.0009586556116119027
$ws.Cells.Item(23, 2).Value = 0.9997809529304504
$ws.Cells.Item(23, 3).Value = 0.000004356679255579365
$ws.Cells.Item(23, 4).Value = 1
$ws.Cells.Item(24, 1).Value = 0.0003222534141968936
$ws.Cells.Item(24, 2).Value = 0.9998783469200134
$ws.Cells.Item(24, 3).Value = 0.000003949335223296657
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(25, 1).Value = 0.000905812019482255
$ws.Cells.Item(25, 2).Value = 0.9998296499252319
$ws.Cells.Item(25, 3).Value = 0.00003876133268931881
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(26, 1).Value = 0.0007430652040056884
$ws.Cells.Item(26, 2).Value = 0.9998539686203003
$ws.Cells.Item(26, 3).Value = 0.00001355701351712923
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(27, 1).Value = 0.0002487579768057913
$ws.Cells.Item(27, 2).Value = 0.9999269843101501
$ws.Cells.Item(27, 3).Value = 0.0000004913151201435539
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(28, 1).Value = 0.001461408683098853
$ws.Cells.Item(28, 2).Value = 0.9997566342353821
$ws.Cells.Item(28, 3).Value = 0.0000004822021537620458
$ws.Cells.Item(28, 4).Value = 1
$ws.Cells.Item(29, 1).Value = 0.000206840006285347
$ws.Cells.Item(29, 2).Value = 0.9999269843101501
$ws.Cells.Item(29, 3).Value = 0.0000009426095175513183
$ws.Cells.Item(29, 4).Value = 1
$ws.Cells.Item(30, 1).Value = 0.0005499968538060784
$ws.Cells.Item(30, 2).Value = 0.9998053312301636
$ws.Cells.Item(30, 3).Value = 0.001929121208377182
$ws.Cells.Item(30, 4).Value = 0.9997029304504395
$ws.Cells.Item(31, 1).Value = 0.0005400101072154939
$ws.Cells.Item(31, 2).Value = 0.9998539686203003
$ws.Cells.Item(31, 3).Value = 0.00000283778695120418
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(32, 1).Value = 0.0003833342925645411
$ws.Cells.Item(32, 2).Value = 0.9999026656150818
$ws.Cells.Item(32, 3).Value = 0.00002175489680666942
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(33, 1).Value = 0.000615773256868124
$ws.Cells.Item(33, 2).Value = 0.9998783469200134
$ws.Cells.Item(33, 3).Value = 0.000004279997938283486
$ws.Cells.Item(33, 4).Value = 1
$ws.Cells.Item(34, 1).Value = 0.0002531058562453836
$ws.Cells.Item(34, 2).Value = 0.9999026656150818
$ws.Cells.Item(34, 3).Value = 0.00000117589911496907
$ws.Cells.Item(34, 4).Value = 1
$ws.Cells.Item(35, 1).Value = 0.000520341913215816
$ws.Cells.Item(35, 2).Value = 0.9999269843101501
$ws.Cells.Item(35, 3).Value = 0.0000003512630257773708
$ws.Cells.Item(35, 4).Value = 1
$ws.Cells.Item(36, 1).Value = 0.0003840986173599958
$ws.Cells.Item(36, 2).Value = 0.9998539686203003
$ws.Cells.Item(36, 3).Value = 0.0000005758133738709148
$ws.Cells.Item(36, 4).Value = 1
$ws.Cells.Item(37, 1).Value = 0.0001344294869340956
$ws.Cells.Item(37, 2).Value = 0.9999513030052185
$ws.Cells.Item(37, 3).Value = 0.0000006366881848407502
$ws.Cells.Item(37, 4).Value = 1
$ws.Cells.Item(38, 1).Value = 0.00008463999984087422
$ws.Cells.Item(38, 2).Value = 0.9999756813049316
$ws.Cells.Item(38, 3).Value = 0.00000003488395350359497
$ws.Cells.Item(38, 4).Value = 1
$ws.Cells.Item(39, 1).Value = 0.0003258000942878425
$ws.Cells.Item(39, 2).Value = 0.9999513030052185
$ws.Cells.Item(39, 3).Value = 0.00000006650922301787432
$ws.Cells.Item(39, 4).Value = 1
$ws.Cells.Item(40, 1).Value = 0.0008395403274334967
$ws.Cells.Item(40, 2).Value = 0.9998783469200134
$ws.Cells.Item(40, 3).Value = 0.0000003503868697407597
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(41, 1).Value = 0.0009779466781765223
$ws.Cells.Item(41, 2).Value = 0.9998783469200134
$ws.Cells.Item(41, 3).Value = 0.0000001967914897704759
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(42, 1).Value = 0.00009595196752343327
$ws.Cells.Item(42, 2).Value = 0.9999756813049316
$ws.Cells.Item(42, 3).Value = 0.0000001359487953322969
$ws.Cells.Item(42, 4).Value = 1
$ws.Cells.Item(43, 1).Value = 0.0001290240325033665
$ws.Cells.Item(43, 2).Value = 0.9999026656150818
$ws.Cells.Item(43, 3).Value = 0.00000009214803498025503
$ws.Cells.Item(43, 4).Value = 1
$ws.Cells.Item(44, 1).Value = 0.0003877757990267128
$ws.Cells.Item(44, 2).Value = 0.9999513030052185
$ws.Cells.Item(44, 3).Value = 0.00000006714410005770333
$ws.Cells.Item(44, 4).Value = 1
$ws.Cells.Item(45, 1).Value = 0.00002724278238019906
$ws.Cells.Item(45, 2).Value = 1
$ws.Cells.Item(45, 3).Value = 0.00000000849974313155144
$ws.Cells.Item(45, 4).Value = 1
$ws.Cells.Item(46, 1).Value = 0.0004350824747234583
$ws.Cells.Item(46, 2).Value = 0.9999269843101501
$ws.Cells.Item(46, 3).Value = 0.000005373336080083391
$ws.Cells.Item(46, 4).Value = 1
$ws.Cells.Item(47, 1).Value = 0.0009376524249091744
$ws.Cells.Item(47, 2).Value = 0.9998296499252319
$ws.Cells.Item(47, 3).Value = 0.00000311507619699114
$ws.Cells.Item(47, 4).Value = 1
$ws.Cells.Item(48, 1).Value = 0.0001457601611036807
$ws.Cells.Item(48, 2).Value = 0.9999269843101501
$ws.Cells.Item(48, 3).Value = 0.000000454628747093011
$ws.Cells.Item(48, 4).Value = 1
$ws.Cells.Item(49, 1).Value = 0.0001182894775411114
$ws.Cells.Item(49, 2).Value = 0.9999756813049316
$ws.Cells.Item(49, 3).Value = 0.000000004249876006667819
$ws.Cells.Item(49, 4).Value = 1
$ws.Cells.Item(50, 1).Value = 0.0004442113277036697
$ws.Cells.Item(50, 2).Value = 0.9998539686203003
$ws.Cells.Item(50, 3).Value = 0.00000005556599447231747
$ws.Cells.Item(50, 4).Value = 1
$ws.Cells.Item(51, 1).Value = 0.00006495912384707481
$ws.Cells.Item(51, 2).Value = 0.9999756813049316
$ws.Cells.Item(51, 3).Value = 0.00000005400802294275309
$ws.Cells.Item(51, 4).Value = 1
